$d = $word.ActiveDocument

# --- 1) Update the 2011 Case Study reference line ---
$d.Content.Find.Execute(
    "Files\\2011 Case Study\\Primary Sources_Policy_Strategies\\2010_national_security_strategy",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Files\\2011 Case Study\\CS1_Primary Sources_Policy_Strategies\\2010 National Security Strategy",
    2) | Out-Null

# --- 2) Update the 2015 Case Study reference line ---
$d.Content.Find.Execute(
    "Files\\2015 Case Study\\Primary Sources_Policy_Strategies\\2015 National Military Strategy CLEAN",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Files\\2015 Case Study\\CS2_Primary Sources_Policy_Strategies\\2015 National Military Strategy",
    2) | Out-Null

# --- 3) Append six new paragraphs after the "Anticipate and adapt..." paragraph ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)

$newParagraphsXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:pStyle w:val="TextBody"/><w:bidi w:val="0"/><w:spacing w:before="113" w:after="113"/>' +
  '<w:ind w:left="113" w:right="113" w:hanging="0"/><w:jc w:val="left"/>' +
  '<w:rPr><w:highlight w:val="lightGray"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr>' +
  '<w:t>Files\\2018 Case Study\\CS3_Primary Sources_Policy_Strategies\\2017 National Security Strategy - § 3 references coded [ 0.22% Coverage]</w:t></w:r></w:p>' +

  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:pStyle w:val="TextBody"/><w:bidi w:val="0"/><w:spacing w:before="113" w:after="113"/>' +
  '<w:ind w:left="113" w:right="113" w:hanging="0"/><w:jc w:val="left"/>' +
  '<w:rPr><w:highlight w:val="lightGray"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr>' +
  '<w:t>Reference 1 - 0.10% Coverage</w:t></w:r></w:p>' +

  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:pStyle w:val="TextBody"/><w:bidi w:val="0"/><w:spacing w:before="0" w:after="0"/>' +
  '<w:jc w:val="left"/><w:rPr/></w:pPr>' +
  '<w:r><w:rPr/>' +
  '<w:t>These competitions require the United States to rethink the policies of the past two decades—policies based on the assumption that engagement with rivals and their inclusion in international institutions and global commerce would turn them into benign actors and trustworthy partners. For the most part, this premise turned out to be false.</w:t></w:r></w:p>' +

  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:pStyle w:val="TextBody"/><w:bidi w:val="0"/><w:spacing w:before="113" w:after="113"/>' +
  '<w:ind w:left="113" w:right="113" w:hanging="0"/><w:jc w:val="left"/>' +
  '<w:rPr><w:highlight w:val="lightGray"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr>' +
  '<w:t>Reference 2 - 0.06% Coverage</w:t></w:r></w:p>' +

  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:pStyle w:val="TextBody"/><w:bidi w:val="0"/><w:spacing w:before="0" w:after="0"/>' +
  '<w:jc w:val="left"/><w:rPr/></w:pPr>' +
  '<w:r><w:rPr/>' +
  '<w:t>The competitions and rivalries facing the United States are not passing trends or momentary problems. They are intertwined, long-term challenges that demand our sustained national att ention and commitment.</w:t></w:r></w:p>' +

  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:pStyle w:val="TextBody"/><w:bidi w:val="0"/><w:spacing w:before="113" w:after="113"/>' +
  '<w:ind w:left="113" w:right="113" w:hanging="0"/><w:jc w:val="left"/>' +
  '<w:rPr><w:highlight w:val="lightGray"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr>' +
  '<w:t>Reference 3 - 0.05% Coverage</w:t></w:r></w:p>' +

  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:pStyle w:val="TextBody"/><w:bidi w:val="0"/><w:spacing w:before="0" w:after="0"/>' +
  '<w:jc w:val="left"/><w:rPr/></w:pPr>' +
  '<w:r><w:rPr/>' +
  '<w:t>But to maintain these advantages, build upon our strengths, and unleash the talents of the American people, we must protect four vital national interests in this competitive world.</w:t></w:r></w:p>'

$insertPoint.InsertXML($newParagraphsXml)
